$wb = $excel.ActiveWorkbook

# --- Add the two new worksheets, at the end of the tab order -----------
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$ws4 = $wb.Worksheets.Add($null, $last)
$ws4.Name = "domoti"

$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$ws5 = $wb.Worksheets.Add($null, $last)
$ws5.Name = "contacts"

# --- Populate "domoti" (sheet4) -----------------------------------------
# Values are written in the same order the author originally typed them
# in, so that the shared-string table comes out in the matching order.
$ws4.Range("A1").Value = 'Thermostat'
$ws4.Range("B2").Value = 'Nest'
$ws4.Range("B5").Value = 'note : Ces sont des thermostats qui fonctionne avec le 24 volt'
$ws4.Range("A27").Value = 'site web intéressant'
$ws4.Range("B28").Value = 'smartthome.com'
$ws4.Range("A8").Value = 'Caméra'
$ws4.Range("B9").Value = 'netatmo'
$ws4.Range("B10").Value = 'nest'
$ws4.Range("B11").Value = 'spypoint'
$ws4.Range("B4").Value = 'caleo'
$ws4.Range("A14").Value = 'Station météo'
$ws4.Range("B29").Value = 'homeremote'
$ws4.Range("A17").Value = 'protocole'
$ws4.Range("B18").Value = 'zigbee'
$ws4.Range("B19").Value = 'zwave'
$ws4.Range("C18").Value = '2.4 gh'
$ws4.Range("C19").Value = '900 mh'
$ws4.Range("C4").Value = 'casaconnect quebecois mais pas recommandé'
$ws4.Range("B21").Value = 'note : privilégié  zigbee'
$ws4.Range("A23").Value = 'produit domotique'
$ws4.Range("A33").Value = 'Boite électrique'
$ws4.Range("B34").Value = 'TED permet d''avoir des infos sur sa consommation '
$ws4.Range("B20").Value = 'insteon'
$ws4.Range("B6").Value = 'utiliser un relais ppur plinthe électrique (30$)'
$ws4.Range("A36").Value = 'application'
$ws4.Range("A37").Value = ' '
$ws4.Range("B37").Value = 'IFTT'
$ws4.Range("B24").Value = 'wemo  protoclole upnp(wallmart) LUMIÈRE PRISE ETC ..'
$ws4.Range("B3").Value = 'ecobee'

# --- Populate "contacts" (sheet5) ---------------------------------------
$ws5.Range("A2").Value = 'Locataire à Claude : Gypse'

# --- Repeated values (reuse existing shared strings, no new <si/>) -----
$ws4.Range("B15").Value = 'netatmo'
$ws4.Range("B30").Value = 'openhab'

# --- View state: selections / active sheet / scroll position -----------
$ws4.Range("B4").Select() | Out-Null
$ws5.Range("A2").Select() | Out-Null
$ws5.Range("A3").Select() | Out-Null
$ws5.Activate() | Out-Null
